$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 390.25
$ws.Range("I2").Value = 456.6
$ws.Range("J2").Value = 279.66666
$ws.Range("K2").Value = 456.6
$ws.Range("L2").Value = 279.66666
$ws.Range("M2").Value = -343.6
$ws.Range("N2").Value = -505.66666

$ws.Range("H62").Value = 7700.625
$ws.Range("I62").Value = 3086.4285
$ws.Range("J62").Value = 40000
$ws.Range("K62").Value = 3086.4285
$ws.Range("L62").Value = 40000
$ws.Range("M62").Value = -2462.4285
$ws.Range("N62").Value = -41248

$ws.Range("H65").Value = 7700.625
$ws.Range("I65").Value = 3086.4285
$ws.Range("J65").Value = 40000
$ws.Range("K65").Value = 15432.1425
$ws.Range("L65").Value = 200000
$ws.Range("M65").Value = -12312.1425
$ws.Range("N65").Value = -206240

$ws.Range("H100").Value = 1264.091
$ws.Range("I100").Value = 1190.5
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1190.5
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -649.5
$ws.Range("N100").Value = -3082

$ws.Range("H105").Value = 36028.25
$ws.Range("J105").Value = 36028.25
$ws.Range("L105").Value = 36028.25
$ws.Range("N105").Value = -43016.25

$ws.Range("H116").Value = 2316.6667
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2633.3333
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2633.3333
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -9517.3333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1895.931
$ws.Range("I45").Value = 1640.0952
$ws.Range("J45").Value = 2567.5
$ws.Range("K45").Value = 1640.0952
$ws.Range("L45").Value = 2567.5
$ws.Range("M45").Value = -1263.0952
$ws.Range("N45").Value = -3321.5

$ws.Range("H63").Value = 3761.2
$ws.Range("I63").Value = 2900
$ws.Range("J63").Value = 3976.5
$ws.Range("K63").Value = 2900
$ws.Range("L63").Value = 3976.5
$ws.Range("M63").Value = -2214
$ws.Range("N63").Value = -5348.5

$ws.Range("H66").Value = 3761.2
$ws.Range("I66").Value = 2900
$ws.Range("J66").Value = 3976.5
$ws.Range("K66").Value = 14500
$ws.Range("L66").Value = 19882.5
$ws.Range("M66").Value = -11068
$ws.Range("N66").Value = -26746.5

$ws.Range("H110").Value = 1968.8572
$ws.Range("I110").Value = 2000
$ws.Range("J110").Value = 1891
$ws.Range("K110").Value = 2000
$ws.Range("L110").Value = 1891
$ws.Range("M110").Value = 45
$ws.Range("N110").Value = -5981

$ws.Range("H119").Value = 35000
$ws.Range("J119").Value = 35000
$ws.Range("L119").Value = 35000
$ws.Range("N119").Value = -44676

$ws.Range("H122").Value = 4466083
$ws.Range("I122").Value = 1683.3334
$ws.Range("J122").Value = 17859282
$ws.Range("K122").Value = 5050.0002
$ws.Range("L122").Value = 53577846
$ws.Range("M122").Value = -2600.0002
$ws.Range("N122").Value = -53582746

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 50000
$ws.Range("J40").Value = 50000
$ws.Range("L40").Value = 50000
$ws.Range("N40").Value = -50530

$ws.Range("H99").Value = 727.53845
$ws.Range("I99").Value = 737.0909
$ws.Range("J99").Value = 675
$ws.Range("K99").Value = 737.0909
$ws.Range("L99").Value = 675
$ws.Range("M99").Value = 760.9091
$ws.Range("N99").Value = -3671

$ws.Range("H107").Value = 1996.6666
$ws.Range("I107").Value = 1002.6667
$ws.Range("J107").Value = 3487.6667
$ws.Range("K107").Value = 1002.6667
$ws.Range("L107").Value = 3487.6667
$ws.Range("M107").Value = 917.3333
$ws.Range("N107").Value = -7327.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3963.7334
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3963.7334
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3963.7334
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4553.7334

$ws.Range("H34").Value = 3963.7334
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3963.7334
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3963.7334
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -4367.7334

$ws.Range("H107").Value = 1072.5454
$ws.Range("I107").Value = 1161.9166
$ws.Range("J107").Value = 965.3
$ws.Range("K107").Value = 1161.9166
$ws.Range("L107").Value = 965.3
$ws.Range("M107").Value = 758.0834
$ws.Range("N107").Value = -4805.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H110").Value = 2672.6667
$ws.Range("I110").Value = 1027
$ws.Range("J110").Value = 2790.2144
$ws.Range("K110").Value = 3081
$ws.Range("L110").Value = 8370.643199999999
$ws.Range("M110").Value = 1009
$ws.Range("N110").Value = -16550.6432

$ws.Range("H131").Value = 1254.1082
$ws.Range("I131").Value = 2011
$ws.Range("J131").Value = 973.7778
$ws.Range("K131").Value = 6033
$ws.Range("L131").Value = 2921.3334
$ws.Range("M131").Value = -993
$ws.Range("N131").Value = -13001.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 538915.6
$ws.Range("I61").Value = 16398.643
$ws.Range("K61").Value = 16398.643
$ws.Range("M61").Value = -16196.643

$ws.Range("H113").Value = 538915.6
$ws.Range("I113").Value = 16398.643
$ws.Range("K113").Value = 16398.643
$ws.Range("M113").Value = -14228.643

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 3514
$ws.Range("I6").Value = 525
$ws.Range("J6").Value = 6503
$ws.Range("K6").Value = 525
$ws.Range("L6").Value = 6503
$ws.Range("M6").Value = -410
$ws.Range("N6").Value = -6733

$ws.Range("H10").Value = 2000
$ws.Range("I10").Value = 2000
$ws.Range("K10").Value = 2000
$ws.Range("M10").Value = -1831

$ws.Range("H17").Value = 4998
$ws.Range("I17").Value = 4998
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 4998
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -4826
$ws.Range("N17").ClearContents()

$ws.Range("H70").Value = 40105
$ws.Range("J70").Value = 40105
$ws.Range("L70").Value = 40105
$ws.Range("N70").Value = -40735

$ws.Range("H73").Value = 40105
$ws.Range("J73").Value = 40105
$ws.Range("L73").Value = 40105
$ws.Range("N73").Value = -42289

$ws.Range("H75").Value = 38130
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 38130
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 38130
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -40002

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H78").Value = 38130
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 38130
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 114390
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -123750

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H93").Value = 69000
$ws.Range("J93").Value = 69000
$ws.Range("L93").Value = 69000
$ws.Range("N93").Value = -73992

$ws.Range("H96").Value = 1069.5
$ws.Range("I96").Value = 923.4
$ws.Range("K96").Value = 923.4
$ws.Range("M96").Value = 449.6

$ws.Range("H100").Value = 1222.9166
$ws.Range("I100").Value = 310.2857
$ws.Range("J100").Value = 2500.6
$ws.Range("K100").Value = 620.5714
$ws.Range("L100").Value = 5001.2
$ws.Range("M100").Value = -79.57140000000004
$ws.Range("N100").Value = -6083.2

$ws.Range("H103").Value = 55202
$ws.Range("J103").Value = 55202
$ws.Range("L103").Value = 55202
$ws.Range("N103").Value = -57546

$ws.Range("H104").Value = 19767.25
$ws.Range("J104").Value = 19767.25
$ws.Range("L104").Value = 19767.25
$ws.Range("N104").Value = -26755.25

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 2647.353
$ws.Range("I107").Value = 608.25
$ws.Range("J107").Value = 7541.2
$ws.Range("K107").Value = 1824.75
$ws.Range("L107").Value = 22623.6
$ws.Range("M107").Value = 95.25
$ws.Range("N107").Value = -26463.6

$ws.Range("H109").Value = 78900
$ws.Range("J109").Value = 78900
$ws.Range("L109").Value = 78900
$ws.Range("N109").Value = -81674

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H122").Value = 1976.0385
$ws.Range("I122").Value = 1308.3636
$ws.Range("J122").Value = 2465.6667
$ws.Range("K122").Value = 3925.0908
$ws.Range("L122").Value = 7397.000100000001
$ws.Range("M122").Value = -1475.0908
$ws.Range("N122").Value = -12297.0001

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
